# ------------------------------------------------------------------
# Re-creates the "merge Anna's final product" commit:
#   * fixes the stray "CMSC " (trailing space) subject code on the
#     Schedule sheet so it matches every other "CMSC" row
#   * adds a new "Subject" column (D) to the Coords sheet mapping
#     each building to the course-subject code primarily taught there
#   * re-selects the Schedule sheet / the newly edited rows, the way
#     the workbook was left after the edit
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$coords   = $wb.Worksheets.Item("Coords")
$capacity = $wb.Worksheets.Item("Capacity")

# ---- 1. Schedule sheet: trailing-space typo fix -------------------
# Row 21 (CMSC 291, "Continued Computer Science for Non-Majors") had
# the subject stored as "CMSC " (trailing space); every other row
# uses "CMSC". Make it consistent.
$schedule.Cells.Item(21, 1).Value = "CMSC"

# ---- 2. Coords sheet: new "Subject" column (D) ---------------------
# Header cell - bold, 14pt, same font as the other two headers but
# with no cell border.
$coordsHeader = $coords.Cells.Item(1, 4)
$coordsHeader.Value = "Subject"
$coordsHeader.Font.Bold = $true
$coordsHeader.Font.Size = 14

# Data rows, building -> primary subject code. The assignment order
# below matters (it controls the order new strings are interned),
# so keep it exactly as written.
$coords.Cells.Item(8, 4).Value  = "CMSC"   # Information Technology
$coords.Cells.Item(6, 4).Value  = "ART"    # Fine Arts
$coords.Cells.Item(4, 4).Value  = "BIOL"   # Biological Sciences
$coords.Cells.Item(5, 4).Value  = "CMPE"   # Engineering
$coords.Cells.Item(3, 4).Value  = "HIST"   # Arts & Humanities
$coords.Cells.Item(9, 4).Value  = "MATH"   # Janet & Walter Sondheim
$coords.Cells.Item(8, 4).Value  = "CMSC"   # (Information Technology, unchanged)
$coords.Cells.Item(7, 4).Value  = "BTEC"   # Interdisciplinary Life S
$coords.Cells.Item(9, 4).Value  = "MATH"   # Janet & Walter Sondheim
$coords.Cells.Item(11, 4).Value = "PYSC"   # Math & Psychology
$coords.Cells.Item(12, 4).Value = "CHEM"   # Meyerhoff Chemistry
$coords.Cells.Item(13, 4).Value = "PHYS"   # Physics
$coords.Cells.Item(14, 4).Value = "POLI"   # Public Policy
$coords.Cells.Item(15, 4).Value = "STAT"   # Sherman Hall
$coords.Cells.Item(1 + 1, 4).Value = "ECON" # Administration (row 2)
$coords.Cells.Item(16, 4).Value = "ENGL"   # University Center
$coords.Cells.Item(2 + 1, 4).Value = "A"    # Albert O. Kuhn Library (row 3)
$coords.Cells.Item(10, 4).Value = "B"      # Lecture Hall 1

$coords.PageSetup.Orientation = 1

# ---- 3. Re-apply the view state the file was left in ---------------
$capacity.Select()
$coords.Select()
$coords.Range("D11").Select() | Out-Null

$schedule.Activate()
$schedule.Range("A21").Select() | Out-Null
